$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.068.72"
$ws.Range("E2").Value = "'  -3.73%  "
$ws.Range("D3").Value = "'2.448.82"
$ws.Range("E3").Value = "'  -3.09%  "
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'309.76"
$ws.Range("E5").Value = "'  +0.46%  "
$ws.Range("D6").Value = "'92.08"
$ws.Range("E6").Value = "'  -8.04%  "
$ws.Range("E7").Value = "'  -2.89%  "
$ws.Range("E8").Value = "'  +0.09%  "
$ws.Range("D9").Value = "'0.491"
$ws.Range("E9").Value = "'  -5.83%  "
$ws.Range("D10").Value = "'32.82"
$ws.Range("E10").Value = "'  -7.90%  "
$ws.Range("D11").Value = "'0.0770"
$ws.Range("E11").Value = "'  -4.04%  "
$ws.Range("E12").Value = "'  -0.89%  "
$ws.Range("D13").Value = "'6.89"
$ws.Range("E13").Value = "'  -6.23%  "
$ws.Range("D14").Value = "'2.829.22"
$ws.Range("E14").Value = "'  -3.00%  "
$ws.Range("D15").Value = "'2.445.24"
$ws.Range("E15").Value = "'  -3.84%  "
$ws.Range("D16").Value = "'14.61"
$ws.Range("E16").Value = "'  -4.28%  "
$ws.Range("D17").Value = "'0.772"
$ws.Range("E17").Value = "'  -4.48%  "
$ws.Range("D18").Value = "'41.027.03"
$ws.Range("E18").Value = "'  -3.78%  "
$ws.Range("D19").Value = "'6.21"
$ws.Range("E19").Value = "'  -7.20%  "
$ws.Range("D20").Value = "'0.0₃0906"
$ws.Range("E20").Value = "'  -4.38%  "
$ws.Range("D21").Value = "'10.98"
$ws.Range("E21").Value = "'  -9.74%  "
$ws.Range("D22").Value = "'67.51"
$ws.Range("E22").Value = "'  -2.80%  "
$ws.Range("D23").Value = "'233.36"
$ws.Range("E23").Value = "'  -3.91%  "
$ws.Range("E24").Value = "'  -4.77%  "
$ws.Range("E25").Value = "'  +0.32%  "
$ws.Range("D26").Value = "'1.87"
$ws.Range("E26").Value = "'  -7.17%  "
$ws.Range("D27").Value = "'23.51"
$ws.Range("E27").Value = "'  -7.53%  "
$ws.Range("D28").Value = "'2.19"
$ws.Range("E28").Value = "'  -5.87%  "
$ws.Range("D29").Value = "'9.46"
$ws.Range("E29").Value = "'  -6.35%  "
$ws.Range("D30").Value = "'35.28"
$ws.Range("E30").Value = "'  -7.89%  "
$ws.Range("D31").Value = "'150.20"
$ws.Range("E31").Value = "'  -4.65%  "
$ws.Range("D32").Value = "'5.40"
$ws.Range("E32").Value = "'  -5.72%  "
$ws.Range("D33").Value = "'2.66"
$ws.Range("E33").Value = "'  -5.47%  "
$ws.Range("D34").Value = "'2.54"
$ws.Range("E34").Value = "'  -3.54%  "
$ws.Range("D35").Value = "'0.0730"
$ws.Range("E35").Value = "'  -6.54%  "
$ws.Range("D36").Value = "'2.94"
$ws.Range("E36").Value = "'  -5.88%  "
$ws.Range("D37").Value = "'16.59"
$ws.Range("E37").Value = "'  -7.61%  "
$ws.Range("D38").Value = "'1.83"
$ws.Range("E38").Value = "'  -6.67%  "
$ws.Range("E39").Value = "'  -3.98%  "
$ws.Range("E40").Value = "'  -8.51%  "
$ws.Range("E41").Value = "'  -2.97%  "
$ws.Range("E42").Value = "'  +0.14%  "
$ws.Range("D43").Value = "'19.51"
$ws.Range("E43").Value = "'  -11.64%  "
$ws.Range("D44").Value = "'1.951.47"
$ws.Range("E44").Value = "'  -2.51%  "
$ws.Range("E45").Value = "'  -6.62%  "
$ws.Range("D46").Value = "'2.97"
$ws.Range("E46").Value = "'  -9.05%  "
$ws.Range("E47").Value = "'  -4.67%  "
$ws.Range("D48").Value = "'69.38"
$ws.Range("E48").Value = "'  -3.47%  "
$ws.Range("D49").Value = "'95.16"
$ws.Range("E49").Value = "'  -5.44%  "
$ws.Range("D50").Value = "'0.174"
$ws.Range("E50").Value = "'  -7.58%  "
$ws.Range("D51").Value = "'72.89"
$ws.Range("E51").Value = "'  -7.81%  "

$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Style = "Normal"
